$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 contact: Dakota Myers / drm1022@sru.edu -> Franz Ferdinand / therockband@gmail.com ---
# (Middle initial "R", address, city, state, zip and phone numbers on row 2 stay the same.)

# The hyperlink object model here doesn't expose a working per-item Delete, and the
# Hyperlink.Address getter comes back empty, so capture each surviving hyperlink's
# mailto target from its cell text (identical to the address for every row in this
# sheet) before touching anything, wipe the whole collection, make the edits, and
# then re-add hyperlinks for every row except D2 (whose link is being removed).
$keepRefs = @("D3", "D4", "D5", "D6", "D8", "D7", "D9", "D10")
$targets = @{}
foreach ($ref in $keepRefs) {
    $targets[$ref] = "mailto:" + $ws.Range($ref).Value()
}

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "Franz"
$ws.Range("B2").Value = "Ferdinand"
$ws.Range("D2").Value = "therockband@gmail.com"

foreach ($ref in $keepRefs) {
    $ws.Hyperlinks.Add($ws.Range($ref), $targets[$ref]) | Out-Null
}

# --- Update the active selection shown when the sheet is opened ---
$ws.Range("M5").Select() | Out-Null
